$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (current: A=#Donor_ID, B=Tissue_ID, C=is_normal,
# D=Sample_ID, E=relative_file_path). This shifts D->E and E->F.
$ws.Columns("D").Insert()

# New header for the inserted column D
$ws.Range("D1").Value = "is_normal_for_donor"

# Only the first data row gets a value in the new column ("Y"); rows 3 and 4 stay blank.
$ws.Range("D2").Value = "Y"

# Set column widths to match target layout (closest achievable values; the
# runtime quantizes stored widths to 1/6-character increments)
$ws.Columns("D").ColumnWidth = 17.5
$ws.Columns("E").ColumnWidth = 9.0

# Update the selected cell to D2, matching the new active selection
$ws.Range("D2").Select()
